# "Add nav bar to ALB" - add a new "alb" column (F) to both offset tables on
# Sheet1, mirroring the existing "z" column (E): a base value in the header
# row of each block (F3 / F9) and, below it, offsets from that base using the
# same (row-C / column-D) deltas already computed for the other stations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column D: bold the values (new font + cell style, same 0.000 numFmt) ---
$ws.Range("D3:D12").Font.Bold = $true

# --- New header cell F2 ("alb") ---
$ws.Range("F2").Value = "alb"

# --- Block 1 (rows 3-6) ---
$ws.Range("F3").Formula = "=0.0932"

$ws.Range("F4:F6").NumberFormat = "0.000"
# E$3+$C4 (and down) get the new mixed-reference form so the pattern can be
# filled right into F; FormulaR1C1 on a multi-cell range resolves per-cell
# relative references correctly (unlike .Formula, which is copied verbatim).
$ws.Range("E4:F4").FormulaR1C1 = "=R3C+RC3"
$ws.Range("E5:F6").FormulaR1C1 = "=R3C+RC3"

# --- Block 2 (rows 9-12) ---
$ws.Range("F9").Formula = "=3.877"

$ws.Range("F10:F12").NumberFormat = "0.000"
$ws.Range("E10:F10").FormulaR1C1 = "=R9C+RC3"
$ws.Range("E11:F12").FormulaR1C1 = "=R9C+RC3"

# --- Selection moves to I5 ---
[void]$ws.Range("I5").Select()
